$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "(" + "the" (gramStart/gramEnd) + " original name?)" -> a single
# run "(the original name?)". A Find/Replace with identical replacement text
# normalises the matched range into one run and drops the proofErr markers.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(the original name?)", $true, $false, $false, $false, `
    $false, $true, 1, $false, "(the original name?)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Implement Material Design? In regards to color and slight
# shadows behind objects" -> split into three runs by inserting new text
# " (or one style) with CSS" after "Implement Material Design". Toggling a
# character-formatting property on the freshly-inserted text (and then
# reverting it) pins a run boundary so the insertion keeps its own <w:r/>
# instead of re-merging with its identically-formatted neighbour.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Implement Material Design") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" (or one style) with CSS")
$rng.Font.Bold = 1
$rng.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: move <w:lastRenderedPageBreak/> from the "Add comments to
# f*cking everything" run up to the "Draw class diagram" run.
# ---------------------------------------------------------------------------

# 3a. Strip the marker from "Add comments to f*cking everything" by
#     retyping its text via a placeholder (a pure same-text reassignment is
#     a no-op and would leave the marker in place).
$rngComments = $d.Content
$rngComments.Find.Execute("Add comments to f*cking everything") | Out-Null
$commentsStart = $rngComments.Start
$rngComments.Text = "ZZZPLACEHOLDERZZZ"
$rngCommentsRestore = $d.Range($commentsStart, $commentsStart + 17)
$rngCommentsRestore.Text = "Add comments to f*cking everything"

# 3b. Insert a replacement run carrying the marker plus the original text
#     right before "Draw class diagram", then delete the now-stale original
#     text run that follows it.
$rngDraw = $d.Content
$rngDraw.Find.Execute("Draw class diagram") | Out-Null
$drawStart = $rngDraw.Start
$drawEnd = $rngDraw.End
$drawLen = $drawEnd - $drawStart

$insPoint = $d.Range($drawStart, $drawStart)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Draw class diagram</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($xml)

$oldDraw = $d.Range($drawStart + $drawLen, $drawEnd + $drawLen)
$oldDraw.Delete()
